# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 836f5434-... file across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 836f5434-...md (row 4)
$wsOverview.Range("G4").Value = "2016-08-12 03:06:59"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 836f5434-... (row 4)
$wsZhCn.Range("H4").Value = "2016-08-12 03:06:54"
$wsZhCn.Range("K4").Value = "2016-08-12 03:07:25"

# de-de sheet: Correspond Handoff Datetime stays the same value as Overview's
# Latest HO Xliff Generate Date (shared string), Correspond Handback DateTime
# for 836f5434-... (row 4) gets updated
$wsDeDe.Range("H4").Value = "2016-08-12 03:06:59"
$wsDeDe.Range("K4").Value = "2016-08-12 03:07:32"
